$d = $word.ActiveDocument

# The commit adds a new, plain paragraph "Filter? By renter" right after
# the paragraph "Admins can view all rental listings and have the option
# to remove inappropriate posts." (everything else in the target diff is
# cosmetic proof-reading-tag / run-merge noise already present in the
# source document).

$anchorText = "Admins can view all rental listings and have the option to remove inappropriate posts."

# Locate the anchor paragraph using Find (robust against exact index).
$find = $d.Content.Find
$found = $find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $anchorEnd = $find.Parent.End

    # Insert a brand new empty paragraph right after the anchor paragraph.
    $insertionPoint = $d.Range($anchorEnd, $anchorEnd)
    $insertionPoint.InsertParagraphAfter()

    # Re-locate the freshly created paragraph through the Paragraphs
    # collection (re-fetching via the collection keeps the object
    # correctly bound, unlike ad-hoc Range(start,end) objects here).
    $paraCount = $d.Paragraphs.Count
    $newParaIndex = -1
    for ($i = 1; $i -le $paraCount; $i++) {
        $candidate = $d.Paragraphs.Item($i).Range
        if ($candidate.Start -eq ($anchorEnd + 1)) {
            $newParaIndex = $i
        }
    }

    if ($newParaIndex -ne -1) {
        # Set the new paragraph's text.
        $d.Paragraphs.Item($newParaIndex).Range.Text = "Filter? By renter"

        # Make sure the new paragraph is plain text (no inherited
        # highlight/bold coming from the preceding paragraph).
        $d.Paragraphs.Item($newParaIndex).Range.HighlightColorIndex = 0
    }
}
